# Override "Transporte" buffer classification = 75 metros
# for the affected station rows, across the 250m/500m/750m/1000m
# buffer columns (H, I, J, K).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 20 (Fernando Pó) -> Outros
$ws.Range("H20:K20").Value = "Outros"

# Rows that become "Urbanizado"
$ws.Range("H29:K29").Value = "Urbanizado"
$ws.Range("H38:K38").Value = "Urbanizado"
$ws.Range("H39:K39").Value = "Urbanizado"
$ws.Range("H45:K45").Value = "Urbanizado"
$ws.Range("H46:K46").Value = "Urbanizado"
